$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell used to stage text-formatted values so the pasted-in value
# lands as literal text (matching the source workbook's inline-string
# cells) instead of being auto-coerced to a number/percentage by Excel's
# normal text-to-number recognition. The helper lives far outside the
# sheet's used range and is cleared (format + contents) at the end so it
# leaves no trace on the sheet's dimensions or styles.
$helper = $ws.Range("ZZ1")
$helper.NumberFormat = "@"

function Set-TextValue([string]$cellRef, [string]$text) {
    $helper.Value = $text
    $helper.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)  # xlPasteValues
}

Set-TextValue "D2" "305.91"
Set-TextValue "E2" "1.01%"
Set-TextValue "D3" "35.82"
Set-TextValue "E3" "1.56%"
Set-TextValue "D4" "5.083"
Set-TextValue "E4" "-0.07%"
Set-TextValue "D5" "0.08077"
Set-TextValue "E5" "0.92%"
Set-TextValue "D6" "1.925"
Set-TextValue "E6" "-0.56%"
Set-TextValue "D7" "4.191"
Set-TextValue "E7" "3.17%"
Set-TextValue "D8" "7.759"
Set-TextValue "E8" "-0.14%"
Set-TextValue "D9" "0.9269"
Set-TextValue "E9" "0.44%"
Set-TextValue "D10" "0.1373"
Set-TextValue "E10" "11.90%"
Set-TextValue "D11" "0.1904"
Set-TextValue "E11" "2.26%"
Set-TextValue "D12" "0.09180"
Set-TextValue "E12" "-3.17%"
Set-TextValue "D13" "0.03420"
Set-TextValue "E13" "-5.96%"
Set-TextValue "D14" "0.09837"
Set-TextValue "E14" "-0.28%"
Set-TextValue "D15" "0.001409"
Set-TextValue "E15" "0.89%"
Set-TextValue "D16" "0.005766"
Set-TextValue "E16" "-0.47%"
Set-TextValue "D17" "3.625"
Set-TextValue "E17" "3.62%"
Set-TextValue "D18" "3.012"
Set-TextValue "E18" "1.80%"
Set-TextValue "D19" "0.3457"
Set-TextValue "E19" "1.54%"
Set-TextValue "D20" "0.1305"
Set-TextValue "E20" "0.22%"
Set-TextValue "E21" "-2.54%"
Set-TextValue "E22" "-0.93%"
Set-TextValue "D23" "0.04439"
Set-TextValue "E23" "-1.96%"
Set-TextValue "D24" "0.001221"
Set-TextValue "E24" "0.38%"
Set-TextValue "D25" "0.004807"
Set-TextValue "E25" "-0.64%"
Set-TextValue "D26" "0.0001431"
Set-TextValue "E26" "14.35%"
Set-TextValue "D27" "0.0003134"
Set-TextValue "E27" "-25.25%"
Set-TextValue "E39" "4.86%"
Set-TextValue "D40" "0.04926"
Set-TextValue "E40" "4.40%"
Set-TextValue "D41" "0.007640"
Set-TextValue "E41" "1.45%"
Set-TextValue "D42" "0.01034"
Set-TextValue "E42" "6.33%"
Set-TextValue "D43" "0.1377"
Set-TextValue "E43" "3.66%"
Set-TextValue "D44" "0.002104"
Set-TextValue "E44" "-0.43%"
Set-TextValue "D45" "0.01107"
Set-TextValue "E45" "9.21%"
Set-TextValue "D46" "0.00006425"
Set-TextValue "E46" "2.28%"
Set-TextValue "D47" "0.00000000751"
Set-TextValue "E47" "0.04%"
Set-TextValue "D48" "64.67"
Set-TextValue "E48" "0.29%"
Set-TextValue "E49" "-19.96%"
Set-TextValue "D50" "0.00002104"
Set-TextValue "E50" "0.04%"
Set-TextValue "D51" "0.0002004"
Set-TextValue "E51" "0.04%"

$helper.Clear()
$excel.CutCopyMode = $false
